$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the (encrypted) message text stored in B3 with a new value.
$ws.Range("B3").Value = "U2FsdGVkX1+OSmLAeytfRFMzwftOg7BmHQT6GYLtwb17ew7I/rZrHXM2KOujqXwc/aBvndfAz6U5ltHv8g67vD7DOQSPPIUTYK2XM/4+XGMNqlxtciFXSfm7nytpBIQOWC3qB+e8ZqUaPUkacfkVQa3MwERnX8LBaNdq3C5bGlqDxO951iCgbOLjn+LuohvCFk0+u4F9shF8B/yPl2YJOZpzSQaWk8hmd7zfHrNOn6aUBWGJi7EJicqJX19/EsyqniFCxBAPKnUs6R3AanXuIXsKGPAwomOLZtwc2yQm6+GHPRacQ0HoT3LEfWC6VjjRWvNWV31EO26KUM3dWj1VRUULSQjlUzOWgYr+8solE0oMz8sxIlda08Bc0ERV4w3t"

# Move the active selection to B6 (matches the saved view state in the edited file).
$ws.Range("B6").Select()
